# atualização do dia 31/10
$wb = $excel.ActiveWorkbook

# --- Sheet "Planejamento" ---
$ws1 = $wb.Worksheets.Item("Planejamento")

# Highlighted period selector
$ws1.Range("G2").Value = 4

# Row 5 - "1. Revisar e atualizar os codigos de SDL"
$ws1.Range("D5").Value = 1
$ws1.Range("E5").Value = 1
$ws1.Range("F5").Value = 1

# Row 6 - "2. Separar o SDL dos simuladores"
$ws1.Range("D6").Value = 3
$ws1.Range("E6").Value = 1
$ws1.Range("F6").Value = 1

# Row 9 - "5. Assets"
$ws1.Range("B9").Value = 3
$ws1.Range("D9").Value = 3
$ws1.Range("F9").Value = 0.2

# --- Sheet "Descrição das Atividades" ---
$ws2 = $wb.Worksheets.Item("Descrição das Atividades")
$ws2.Range("B6:P6").Select()

# Restore "Planejamento" as the active sheet/selection
$ws1.Select()
$ws1.Range("E5").Select()
